# Disaggregation of commodity Copper
# Rename "Copper ores and concentrates" -> "Copper" in cell C4 of every
# yearly worksheet, and update the handful of D4 figures whose underlying
# source values shifted by a tiny (last-significant-digit) amount as part
# of the same disaggregation update.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("C4")
    if ($cell.Value() -eq "Copper ores and concentrates") {
        $cell.Value = "Copper"
    }
}

$years  = @("2021", "2023", "2025", "2028", "2031", "2041", "2044", "2048", "2054", "2072", "2077", "2092")
$values = @(27629.08234046596, 45474.52846901826, 52615.60445701829, 76414.31688861702, 102815.2737995718, 413896.0450908013, 824812.1633242127, 1595482.848064659, 1715249.128188553, 1598520.870762428, 1443658.30372603, 1706284.654525028)

for ($i = 0; $i -lt $years.Length; $i++) {
    $ws = $wb.Worksheets.Item($years[$i])
    $ws.Range("D4").Value = $values[$i]
}
